$wb = $excel.ActiveWorkbook

# Column order matches header row: eb, hp, st, wi, gt, dgt, ieh, chp, ac, ab_ct, ab_hp, cp_ct, cp_hp, ttes, ites
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")

# New values for row 2, columns A..O, one array per sheet (order matches $sheetNames)
$sheetNames = @("2025","2030","2035","2040","2045","2050")

$values2025 = @(0, 10372.65132737054, 0, 0, 289260.5393052954, 0, 80959.25712661834, 0, 161710.6685703679, 0, 0, 484922.2142001599, 105953.7713982, 70003.73489578845, 69744.89343456978)
$values2030 = @(0, 31203.23858116339, 0, 0, 170658.5511254234, 0, 0, 0, 209080.6134235085, 0, 0, 63518.11613148725, 68536.72857011756, 19285.19160463996, 27033.1386905727)
$values2035 = @(27543.1755456332, 22113.21643273498, 0, 0, 114655.4402706629, 0, 0, 0, 153866.0861464091, 0, 0, 0, 44638.22942194272, 39676.88529639924, 31311.04369977792)
$values2040 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1142.580190039942, 0)
$values2045 = @(29588.33508286276, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 4347.543515635315, 20429.76977394434)
$values2050 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

$allValues = @($values2025, $values2030, $values2035, $values2040, $values2045, $values2050)

for ($s = 0; $s -lt $sheetNames.Length; $s++) {
    $ws = $wb.Worksheets.Item([string]$sheetNames[$s])
    $rowValues = $allValues[$s]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + "2").Value = $rowValues[$i]
    }
}
